# Auto-generated edit script: refresh market-price derived cells across all 8 class sheets
# (values sourced from scheduled runner data pull; no formulas involved, plain overwrite)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K38").Value = 552
$ws.Range("H38").Value = 3110
$ws.Range("I38").Value = 184
$ws.Range("M38").Value = -180
$ws.Range("J74").Value = 1300
$ws.Range("H74").Value = 3223.5334
$ws.Range("L74").Value = 1300
$ws.Range("I74").Value = 3360.9285
$ws.Range("K74").Value = 3360.9285
$ws.Range("N74").Value = -3172
$ws.Range("M74").Value = -2424.9285
$ws.Range("L77").Value = 6500
$ws.Range("K77").Value = 16804.6425
$ws.Range("I77").Value = 3360.9285
$ws.Range("H77").Value = 3223.5334
$ws.Range("J77").Value = 1300
$ws.Range("N77").Value = -15860
$ws.Range("M77").Value = -12124.6425
$ws.Range("L94").Value = 34998
$ws.Range("H94").Value = 21054.223
$ws.Range("N94").Value = -35900
$ws.Range("I94").Value = 3624.5
$ws.Range("K94").Value = 3624.5
$ws.Range("J94").Value = 34998
$ws.Range("M94").Value = -3173.5
$ws.Range("N112").Value = -9338
$ws.Range("H112").Value = 2075.0557
$ws.Range("J112").Value = 2374
$ws.Range("L112").Value = 7122
$ws.Range("K132").Value = 4473.72
$ws.Range("I132").Value = 1491.24
$ws.Range("H132").Value = 358598.4
$ws.Range("M132").Value = -1943.72
$ws.Range("I137").Value = 2163.7273
$ws.Range("H137").Value = 4479.636
$ws.Range("K137").Value = 6491.1819
$ws.Range("M137").Value = -3941.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1049.5217
$ws.Range("M2").Value = -936.5217
$ws.Range("K2").Value = 1049.5217
$ws.Range("H2").Value = 1072.4584
$ws.Range("J32").Value = 59526.816
$ws.Range("M32").Value = -2476.1396
$ws.Range("K32").Value = 2763.1396
$ws.Range("L32").Value = 59526.816
$ws.Range("N32").Value = -60100.816
$ws.Range("H32").Value = 14326.111
$ws.Range("I32").Value = 2763.1396
$ws.Range("N105").Value = -95651
$ws.Range("J105").Value = 88663
$ws.Range("H105").Value = 88663
$ws.Range("L105").Value = 88663
$ws.Range("I110").Value = 3001
$ws.Range("H110").Value = 3299.8333
$ws.Range("K110").Value = 3001
$ws.Range("M110").Value = -956
$ws.Range("H116").Value = 1072.4584
$ws.Range("M116").Value = 1244.4783
$ws.Range("K116").Value = 1049.5217
$ws.Range("I116").Value = 1049.5217
$ws.Range("M122").Value = -5424.099999999999
$ws.Range("K122").Value = 7874.099999999999
$ws.Range("H122").Value = 2762.2964
$ws.Range("I122").Value = 2624.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1072.4584
$ws.Range("K3").Value = 1049.5217
$ws.Range("I3").Value = 1049.5217
$ws.Range("M3").Value = -935.5217
$ws.Range("H26").Value = 20117.25
$ws.Range("J26").Value = 30000
$ws.Range("L26").Value = 30000
$ws.Range("I26").Value = 16823
$ws.Range("K26").Value = 16823
$ws.Range("N26").Value = -30584
$ws.Range("M26").Value = -16531
$ws.Range("H64").Value = 2033
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 2033
$ws.Range("L67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H107").Value = 1466.2
$ws.Range("K107").Value = 1441.2778
$ws.Range("I107").Value = 1441.2778
$ws.Range("M107").Value = 478.7221999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10833.333
$ws.Range("K7").Value = 142912
$ws.Range("M7").Value = -142799
$ws.Range("H7").Value = 62915.312
$ws.Range("I7").Value = 142912
$ws.Range("M31").Value = -3505.8
$ws.Range("H31").Value = 4827.316
$ws.Range("I31").Value = 3800.8
$ws.Range("K31").Value = 3800.8
$ws.Range("K34").Value = 3800.8
$ws.Range("I34").Value = 3800.8
$ws.Range("M34").Value = -3598.8
$ws.Range("H34").Value = 4827.316
$ws.Range("N107").Value = -5039.75
$ws.Range("H107").Value = 1050.8334
$ws.Range("K107").Value = 976.375
$ws.Range("J107").Value = 1199.75
$ws.Range("I107").Value = 976.375
$ws.Range("L107").Value = 1199.75
$ws.Range("M107").Value = 943.625
$ws.Range("M134").Value = -8864.6844
$ws.Range("K134").Value = 11399.6844
$ws.Range("I134").Value = 3799.8948
$ws.Range("H134").Value = 4628.476

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K9").Value = 6298.5
$ws.Range("J9").Value = 32208.455
$ws.Range("N9").Value = -97073.36500000001
$ws.Range("M9").Value = -6074.5
$ws.Range("H9").Value = 24179.4
$ws.Range("L9").Value = 96625.36500000001
$ws.Range("I9").Value = 2099.5
$ws.Range("N17").Value = -4800.5
$ws.Range("H17").Value = 1322.3334
$ws.Range("L17").Value = 4462.5
$ws.Range("J17").Value = 1487.5
$ws.Range("L34").Value = 1797.9999
$ws.Range("J34").Value = 599.3333
$ws.Range("K34").Value = 1344.42855
$ws.Range("I34").Value = 448.14285
$ws.Range("M34").Value = -1260.42855
$ws.Range("N34").Value = -1965.9999
$ws.Range("H34").Value = 474.82352
$ws.Range("J39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M99").Value = -754
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 3000
$ws.Range("L114").Value = 3747.75
$ws.Range("H114").Value = 1582.1666
$ws.Range("N114").Value = -10255.75
$ws.Range("J114").Value = 1249.25
$ws.Range("K132").Value = 125984.997
$ws.Range("I132").Value = 13998.333
$ws.Range("H132").Value = 10052.667
$ws.Range("M132").Value = -123454.997
$ws.Range("K139").Value = 13190.0625
$ws.Range("I139").Value = 4396.6875
$ws.Range("M139").Value = -8050.0625
$ws.Range("H139").Value = 8737.130999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I43").Value = 3456.125
$ws.Range("H43").Value = 6405.4443
$ws.Range("K43").Value = 3456.125
$ws.Range("M43").Value = -3305.125
$ws.Range("J48").Value = 48500
$ws.Range("L48").Value = 48500
$ws.Range("N48").Value = -49470
$ws.Range("H48").Value = 35666.668
$ws.Range("J52").Value = 57450
$ws.Range("N52").Value = -57968
$ws.Range("L52").Value = 57450
$ws.Range("H52").Value = 54310
$ws.Range("I102").Value = 43401.477
$ws.Range("J102").Value = 8998.666999999999
$ws.Range("M102").Value = -41779.477
$ws.Range("N102").Value = -12242.667
$ws.Range("K102").Value = 43401.477
$ws.Range("L102").Value = 8998.666999999999
$ws.Range("H102").Value = 39101.125
$ws.Range("H107").Value = 1107.2858
$ws.Range("K107").Value = 758.44446
$ws.Range("I107").Value = 758.44446
$ws.Range("M107").Value = 1161.55554
$ws.Range("K113").Value = 433
$ws.Range("H113").Value = 433
$ws.Range("M113").Value = 1737
$ws.Range("I113").Value = 433
$ws.Range("J123").Value = 50326
$ws.Range("L123").Value = 50326
$ws.Range("N123").Value = -55226
$ws.Range("H123").Value = 37807.25
$ws.Range("N126").Value = -18431.375
$ws.Range("M126").Value = -11103.6362
$ws.Range("H126").Value = 4513
$ws.Range("L126").Value = 13491.375
$ws.Range("I126").Value = 4524.5454
$ws.Range("J126").Value = 4497.125
$ws.Range("K126").Value = 13573.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 431.72726
$ws.Range("M16").Value = -261.72726
$ws.Range("K16").Value = 431.72726
$ws.Range("I16").Value = 431.72726
$ws.Range("N88").Value = -35854.5
$ws.Range("L88").Value = 34998.5
$ws.Range("J88").Value = 34998.5
$ws.Range("H88").Value = 28749
$ws.Range("N91").Value = -37962.5
$ws.Range("H91").Value = 28749
$ws.Range("J91").Value = 34998.5
$ws.Range("L91").Value = 34998.5
$ws.Range("H100").Value = 3865.6667
$ws.Range("M100").Value = -3124.6667
$ws.Range("K100").Value = 3665.6667
$ws.Range("I100").Value = 3665.6667
$ws.Range("J100").Value = 4065.6667
$ws.Range("N100").Value = -5147.6667
$ws.Range("L100").Value = 4065.6667
$ws.Range("K132").Value = 7200
$ws.Range("I132").Value = 2400
$ws.Range("H132").Value = 2626.9092
$ws.Range("M132").Value = -4670

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("K2").Value = 0
$ws.Range("N2").Value = -19557
$ws.Range("J2").Value = 19333
$ws.Range("H2").Value = 19333
$ws.Range("L2").Value = 19333
$ws.Range("L4").Value = 2100.5
$ws.Range("K4").Value = 4173
$ws.Range("J4").Value = 2100.5
$ws.Range("N4").Value = -2326.5
$ws.Range("H4").Value = 2791.3333
$ws.Range("I4").Value = 4173
$ws.Range("M4").Value = -4060
$ws.Range("J29").Value = 44000
$ws.Range("N29").Value = -44580
$ws.Range("H29").Value = 45666.668
$ws.Range("L29").Value = 44000
$ws.Range("N107").Value = -10444.2855
$ws.Range("H107").Value = 1516.3572
$ws.Range("K107").Value = 2493.8571
$ws.Range("J107").Value = 2201.4285
$ws.Range("I107").Value = 831.2857
$ws.Range("L107").Value = 6604.2855
$ws.Range("M107").Value = -573.8571000000002
$ws.Range("K113").Value = 1092
$ws.Range("H113").Value = 719.36365
$ws.Range("M113").Value = 1078
$ws.Range("I113").Value = 364
